$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of data (2-10) have been re-shuffled: new row N gets the
# D, L, M, N, O, P, Q, R, S, T values that used to live in old row Map[N].
# Columns A, B, C, E, F, G, H, I, J, K stay identical (they were the same
# for every row already).

$rowMap = @{
    2  = 4
    3  = 9
    4  = 10
    5  = 3
    6  = 7
    7  = 5
    8  = 6
    9  = 8
    10 = 2
}

# Snapshot the "before" values of the columns that move, keyed by source row.
$colIdx = @{ D = 4; L = 12; M = 13; N = 14; O = 15; P = 16; Q = 17; R = 18; S = 19; T = 20 }

$snapshot = @{}
foreach ($r in 2..10) {
    $rowData = @{}
    foreach ($colName in $colIdx.Keys) {
        $c = $colIdx[$colName]
        $rowData[$colName] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

foreach ($destRow in 2..10) {
    $srcRow = $rowMap[$destRow]
    $src = $snapshot[$srcRow]
    foreach ($colName in $colIdx.Keys) {
        $c = $colIdx[$colName]
        $ws.Cells.Item($destRow, $c).Value2 = $src[$colName]
    }
}
